$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# New header cell O1: "Serviced by", copying the format of the adjacent
# header cell (N1, "Correction") so it picks up the same bold/bordered style.
$ws.Cells.Item(1, 15).Value = "Serviced by"
$ws.Cells.Item(1, 14).Copy()
$ws.Cells.Item(1, 15).PasteSpecial(-4122)

# New (empty) column cells O2:O12 for the existing data rows, present but blank.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 15).Style = "Normal"
}

# The old trailing fully-blank row 13 is removed entirely.
$ws.Rows.Item(13).Delete()
